$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 860.1177
$ws.Range("I15").Value = 860.1177
$ws.Range("K15").Value = 2580.3531
$ws.Range("M15").Value = -2411.3531
$ws.Range("H33").Value = 245
$ws.Range("I33").Value = 194
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 194
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 35
$ws.Range("N33").Value = -958
$ws.Range("H101").Value = 1185.2858
$ws.Range("I101").Value = 1259.6
$ws.Range("K101").Value = 3778.8
$ws.Range("M101").Value = -2156.8
$ws.Range("H116").Value = 4299.2
$ws.Range("I116").Value = 4499
$ws.Range("K116").Value = 4499
$ws.Range("M116").Value = -1057
$ws.Range("H131").Value = 2449.5833
$ws.Range("I131").Value = 2089.5
$ws.Range("K131").Value = 6268.5
$ws.Range("M131").Value = -1228.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2416.6
$ws.Range("I2").Value = 1796.125
$ws.Range("K2").Value = 1796.125
$ws.Range("M2").Value = -1683.125
$ws.Range("H32").Value = 6559.2104
$ws.Range("I32").Value = 3570.5312
$ws.Range("K32").Value = 3570.5312
$ws.Range("M32").Value = -3283.5312
$ws.Range("H45").Value = 1998.3636
$ws.Range("I45").Value = 1998.3
$ws.Range("K45").Value = 1998.3
$ws.Range("M45").Value = -1621.3
$ws.Range("H74").Value = 1183.1818
$ws.Range("I74").Value = 1066.7778
$ws.Range("J74").Value = 1707
$ws.Range("K74").Value = 1066.7778
$ws.Range("L74").Value = 1707
$ws.Range("M74").Value = -192.7778000000001
$ws.Range("N74").Value = -3455
$ws.Range("H77").Value = 1183.1818
$ws.Range("I77").Value = 1066.7778
$ws.Range("J77").Value = 1707
$ws.Range("K77").Value = 5333.889
$ws.Range("L77").Value = 8535
$ws.Range("M77").Value = -965.8890000000001
$ws.Range("N77").Value = -17271
$ws.Range("H116").Value = 2416.6
$ws.Range("I116").Value = 1796.125
$ws.Range("K116").Value = 1796.125
$ws.Range("M116").Value = 497.875
$ws.Range("H122").Value = 2252.25
$ws.Range("I122").Value = 2086.3333
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 6258.999899999999
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -3808.999899999999
$ws.Range("N122").Value = -13150
$ws.Range("H135").Value = 37500
$ws.Range("J135").Value = 37500
$ws.Range("L135").Value = 37500
$ws.Range("N135").Value = -47640
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2416.6
$ws.Range("I3").Value = 1796.125
$ws.Range("K3").Value = 1796.125
$ws.Range("M3").Value = -1682.125
$ws.Range("H99").Value = 34916.234
$ws.Range("I99").Value = 44942.957
$ws.Range("K99").Value = 44942.957
$ws.Range("M99").Value = -43444.957

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("M4").Value = 111
$ws.Range("H31").Value = 4636.364
$ws.Range("I31").Value = 3885.625
$ws.Range("K31").Value = 3885.625
$ws.Range("M31").Value = -3590.625
$ws.Range("H34").Value = 4636.364
$ws.Range("I34").Value = 3885.625
$ws.Range("K34").Value = 3885.625
$ws.Range("M34").Value = -3683.625
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -112488
$ws.Range("H99").Value = 11531.481
$ws.Range("I99").Value = 8403.375
$ws.Range("K99").Value = 8403.375
$ws.Range("M99").Value = -6905.375
$ws.Range("H126").Value = 11531.481
$ws.Range("I126").Value = 8403.375
$ws.Range("K126").Value = 25210.125
$ws.Range("M126").Value = -22740.125
$ws.Range("H134").Value = 2102.4
$ws.Range("I134").Value = 1170.3334
$ws.Range("K134").Value = 3511.0002
$ws.Range("M134").Value = -976.0001999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 299.89474
$ws.Range("I12").Value = 294.6111
$ws.Range("J12").Value = 304.65
$ws.Range("K12").Value = 883.8333
$ws.Range("L12").Value = 913.9499999999999
$ws.Range("M12").Value = -710.8333
$ws.Range("N12").Value = -1259.95
$ws.Range("H140").Value = 1631.8
$ws.Range("I140").Value = 1631.8
$ws.Range("K140").Value = 4895.4
$ws.Range("M140").Value = 284.6000000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 10000
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H132").Value = 2723.4707
$ws.Range("I132").Value = 2230.7
$ws.Range("J132").Value = 3427.4285
$ws.Range("K132").Value = 6692.099999999999
$ws.Range("L132").Value = 10282.2855
$ws.Range("M132").Value = -4162.099999999999
$ws.Range("N132").Value = -15342.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4007800
$ws.Range("J2").Value = 8666.666999999999
$ws.Range("L2").Value = 8666.666999999999
$ws.Range("N2").Value = -8890.666999999999
$ws.Range("H61").Value = 1118.8572
$ws.Range("I61").Value = 1055.3334
$ws.Range("K61").Value = 1055.3334
$ws.Range("M61").Value = -853.3334
$ws.Range("H62").Value = 54210.332
$ws.Range("J62").Value = 52631
$ws.Range("L62").Value = 52631
$ws.Range("N62").Value = -53879
$ws.Range("H65").Value = 54210.332
$ws.Range("J65").Value = 52631
$ws.Range("L65").Value = 157893
$ws.Range("N65").Value = -164133
$ws.Range("H68").Value = 7128.5
$ws.Range("J68").Value = 7375.75
$ws.Range("L68").Value = 7375.75
$ws.Range("N68").Value = -8873.75
$ws.Range("H71").Value = 7128.5
$ws.Range("J71").Value = 7375.75
$ws.Range("L71").Value = 36878.75
$ws.Range("N71").Value = -44366.75
$ws.Range("H76").Value = 34500
$ws.Range("J76").Value = 34500
$ws.Range("L76").Value = 34500
$ws.Range("N76").Value = -35176
$ws.Range("H79").Value = 34500
$ws.Range("J79").Value = 34500
$ws.Range("L79").Value = 34500
$ws.Range("N79").Value = -36840
$ws.Range("H113").Value = 1118.8572
$ws.Range("I113").Value = 1055.3334
$ws.Range("K113").Value = 1055.3334
$ws.Range("M113").Value = 1114.6666
$ws.Range("H122").Value = 5816.3335
$ws.Range("I122").Value = 5680
$ws.Range("J122").Value = 6498
$ws.Range("K122").Value = 17040
$ws.Range("L122").Value = 19494
$ws.Range("M122").Value = -14590
$ws.Range("N122").Value = -24394
$ws.Range("H136").Value = 2449.3125
$ws.Range("I136").Value = 2085
$ws.Range("K136").Value = 6255
$ws.Range("M136").Value = -3705

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 32499.5
$ws.Range("J68").Value = 32499.5
$ws.Range("L68").Value = 32499.5
$ws.Range("N68").Value = -34121.5
$ws.Range("H69").Value = 19490.334
$ws.Range("J69").Value = 25110.5
$ws.Range("L69").Value = 25110.5
$ws.Range("N69").Value = -26608.5
$ws.Range("H71").Value = 32499.5
$ws.Range("J71").Value = 32499.5
$ws.Range("L71").Value = 97498.5
$ws.Range("N71").Value = -105610.5
$ws.Range("H72").Value = 19490.334
$ws.Range("J72").Value = 25110.5
$ws.Range("L72").Value = 75331.5
$ws.Range("N72").Value = -82819.5
$ws.Range("H132").Value = 46147.19
$ws.Range("I132").Value = 67971.42999999999
$ws.Range("K132").Value = 203914.29
$ws.Range("M132").Value = -201384.29
